$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# New monthly values for 2022 (row 2): Outubro (J2) and Novembro (K2)
$ws.Range("J2").Value = 0.6814
$ws.Range("K2").Value = 0.6501

# Shift the "Fonte" (source) footnotes down by one row and add the new
# most-recent source link at the top (row 2), pushing the older ones down.
$ws.Range("O4").Value = $ws.Range("O3").Value()
$ws.Range("O3").Value = $ws.Range("O2").Value()
$ws.Range("O2").Value = "http://www.yahii.com.br/poupanca.html"
